$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update project/benefit description texts: "Pedra de Ferro" -> "Ferro" ---
$ws.Range("G2").Value = "O centro administrativo do acampamento onde são produzidos recursos. Produz 400 Utar, 75 Madeira e 75 Ferro toda semana. Adiciona 3 Segurança, 1 Prosperidade e 9 Torres de Turret. Pode-se recrutar 4 equipes de Servos Ifrit a cada semana. Limite de construção: 2."
$ws.Range("G3").Value = "O coração administrativo do acampamento onde os recursos são produzidos. Produz 600 Utar, 100 Madeira, 100 Ferro e 5 Jades toda semana. Adiciona 4 de Segurança, 2 de Prosperidade e 12 Torres de Turret. É possível recrutar uma equipe de Guerreiro Ifrit, Atirador Ifrit, Lançador de Feitiços Ifrit e Cavaleiro Ifrit toda semana. Limite de construção: 3."
$ws.Range("G4").Value = "O centro administrativo do acampamento onde os recursos são produzidos. Produz 400 Utar, 75 Madeira, 75 Ferro e 75 Ferro toda semana. Adiciona 3 de Segurança, 1 de Prosperidade e 9 Torres de Turret. É possível recrutar 2 equipes de Milícia Nasir e 2 equipes de Infantaria Nasir toda semana. Limite de construção: 2."
$ws.Range("G5").Value = "O centro administrativo do acampamento onde os recursos são produzidos. Produz 600 Utar, 100 Madeira, 100 Ferro e 5 Jades toda semana. Adiciona 4 de Segurança, 2 de Prosperidade e 12 Torres de Turret. É possível recrutar uma equipe de Galante Nasir, Lança Longa Nasir, Cavaleiro Nasir, Atirador Nasir e Curandeiro Nasir toda semana. Limite de construção: 3."
$ws.Range("G6").Value = "O centro administrativo do acampamento onde os recursos são produzidos. Produz 400 Utar, 75 Madeira, 75 Ferro e 75 Ferro toda semana. Adiciona 3 de Segurança, 1 de Prosperidade e 9 Torres de Turret. É possível recrutar 2 equipes de Milícia Nasir e 2 equipes de Infantaria Nasir toda semana. Limite de construção: 2."
$ws.Range("G7").Value = "O coração administrativo do acampamento onde os recursos são produzidos. Produz 600 Utar, 100 Madeira, 100 Ferro e 5 Jades toda semana. Adiciona 4 de Segurança, 2 de Prosperidade e 12 Torres de Turret. É possível recrutar uma equipe de Galante Nasir, Lança Longa Nasir, Cavaleiro Nasir, Atirador Nasir e Curandeiro Nasir toda semana. Limite de construção: 3."
$ws.Range("G8").Value = "O centro administrativo do acampamento onde os recursos são produzidos. Produz 400 Utar, 75 Madeira e 75 Ferro toda semana. Adiciona 3 de Segurança, 1 de Prosperidade e 9 Torres de Turret. Pode-se recrutar 1 equipe de Guarda Dhib, 1 equipe de Caçador Dhib e 2 equipes de Pastores Dhib a cada semana. Limite de construção: 2."
$ws.Range("G9").Value = "O centro administrativo do acampamento onde os recursos são produzidos. Produz 600 Utar, 100 Madeira, 100 Ferro e 5 Jades toda semana. Adiciona 4 de Segurança, 2 de Prosperidade e 12 Torres de Turret. Pode-se recrutar 1 equipe de Espadachim Dhib, 1 equipe de Lança Longa Dhib e 2 equipes de Rastreadores Dhib a cada semana. Limite de construção: 3."
$ws.Range("G10").Value = "O centro administrativo do acampamento onde os recursos são produzidos. Produz 700 Utar, 125 Madeira, 125 Ferro e 8 Jades toda semana. Adiciona 4 de Segurança, 3 de Prosperidade e 12 Torres de Turret. Pode-se recrutar 2 equipes de Espadachins Dhib, 2 equipes de Lança Longa Dhib e 2 equipes de Rastreadores Dhib a cada semana. Limite de construção: 3."
$ws.Range("G11").Value = "O centro administrativo do acampamento onde os recursos são produzidos. Produz 400 Utar, 75 Madeira e 75 Ferro toda semana. Adiciona 3 de Segurança, 1 de Prosperidade e 9 Torres de Turret. É possível recrutar 1 equipe de Espião Dakn e 2 equipes de Lançadores Dakn a cada semana. Limite de construção: 2."
$ws.Range("G12").Value = "O centro administrativo do acampamento onde os recursos são produzidos. Produz 600 Utar, 100 Madeira, 100 Ferro e 5 Jades toda semana. Adiciona 4 de Segurança, 2 de Prosperidade e 12 Torres de Turret. É possível recrutar uma equipe de Assassino Dakn, Lançador de Dardos Dakn, Atirador de Espinhos Dakn e Alquimista Dakn toda semana. Limite de construção: 3."
$ws.Range("G13").Value = "O centro administrativo do acampamento onde os recursos são produzidos. Produz 400 Utar, 75 Madeira, 75 Ferro e 75 Ferro toda semana. Adiciona 3 de Segurança, 1 de Prosperidade e 9 Torres de Turret. É possível recrutar 1 equipe de Guerreiro Thur, 1 equipe de Atirador Thur e 2 equipes de Camponeses Thur a cada semana. Limite de construção: 2."
$ws.Range("G14").Value = "O centro administrativo do acampamento onde os recursos são produzidos. Produz 600 Utar, 100 Madeira, 100 Ferro e 5 Jades toda semana. Adiciona 4 de Segurança, 2 de Prosperidade e 12 Torres de Turret. É possível recrutar 1 equipe de Militante Thur, 1 equipe de Arqueiro Exímio Thur e 2 equipes de Infantaria Pesada Thur a cada semana. Limite de construção: 3."
$ws.Range("G17").Value = "Um posto avançado com capacidades defensivas limitadas e espaço para tropas guarnecidas. Produz 300 Utar, 75 Madeira e 75 Ferro toda semana. Adiciona 3 de Segurança e 9 Torres de Turret. Pode-se recrutar 3 equipes de Milícia Nasir a cada semana. Limite de construção: 1."
$ws.Range("G18").Value = "Um posto avançado com capacidades defensivas limitadas e espaço para tropas guarnecidas. Produz 300 Utar, 75 Madeira e 75 Ferro toda semana. Adiciona 3 de Segurança e 12 Torres de Turret. Pode-se recrutar 3 equipes de Milícia Nasir a cada semana. Limite de construção: 1."
$ws.Range("G19").Value = "Um posto avançado com capacidades defensivas limitadas e espaço para tropas guarnecidas. Produz 300 Utar, 75 Madeira e 75 Ferro toda semana. Adiciona 3 de Segurança e 9 Torres de Turret. Pode-se recrutar 3 equipes de Moradores Akhal a cada semana. Limite de construção: 1."
$ws.Range("G20").Value = "Um posto avançado com capacidades defensivas limitadas e espaço para tropas guarnecidas. Produz 300 Utar, 75 Madeira e 75 Ferro toda semana. Adiciona 3 de Segurança e 12 Torres de Turret. Pode-se recrutar 1 equipe de Cavalaria Leve Akhal e 3 equipes de Moradores Akhal a cada semana. Limite de construção: 1."
$ws.Range("G21").Value = "Um posto avançado com capacidades defensivas limitadas e espaço para tropas guarnecidas. Produz 300 Utar, 75 Madeira e 75 Ferro toda semana. Adiciona 3 de Segurança e 9 Torres de Turret. Pode-se recrutar 3 equipes de Pastores Dhib a cada semana. Limite de construção: 1."
$ws.Range("G22").Value = "Um posto avançado com capacidades defensivas limitadas e espaço para tropas guarnecidas. Produz 300 Utar, 75 Madeira, 75 Ferro e 1 Jade toda semana. Adiciona 3 de Segurança e 12 Torres de Turret. Pode-se recrutar 3 equipes de Pastores Dhib a cada semana. Limite de construção: 1."
$ws.Range("G23").Value = "Um posto avançado com capacidades defensivas limitadas e espaço para tropas guarnecidas. Produz 300 Utar, 75 Madeira e 75 Ferro toda semana. Adiciona 3 de Segurança e 9 Torres de Turret. Pode-se recrutar 3 equipes de Camponeses Thur a cada semana. Limite de construção: 1."
$ws.Range("G24").Value = "Um posto avançado com capacidades defensivas limitadas e espaço para tropas guarnecidas. Produz 300 Utar, 75 Madeira, 75 Ferro e 1 Jade toda semana. Adiciona 3 de Segurança e 12 Torres de Turret. Pode-se recrutar 3 equipes de Camponeses Thur a cada semana. As equipes guarnecidas aqui recebem 30 EXP por dia. Limite de construção: 1."
$ws.Range("G25").Value = "Um posto avançado com capacidades defensivas limitadas e espaço para tropas guarnecidas. Produz 300 Utar, 75 Madeira e 75 Ferro toda semana. Adiciona 3 de Segurança e 9 Torres de Turret. Pode-se recrutar 1 equipe de Lançador Dakn e 2 equipes de Donzelas Dakn a cada semana. Limite de construção: 1."
$ws.Range("G26").Value = "Um posto avançado com capacidades defensivas limitadas e espaço para tropas guarnecidas. Produz 300 Utar, 75 Madeira e 75 Ferro toda semana. Adiciona 3 de Segurança e 12 Torres de Turret. Pode-se recrutar 1 equipe de Lançador Dakn, 1 equipe de Escorpião Tóxico e 2 equipes de Donzelas Dakn a cada semana. Limite de construção: 1."
$ws.Range("G29").Value = "Produz 20 unidades de Ferro toda semana. "
$ws.Range("G30").Value = "Produz 50 unidades de Ferro toda semana. "
$ws.Range("G37").Value = "Produz 20 unidades de Ferro toda semana. "
$ws.Range("G38").Value = "Produz 50 unidades de Ferro toda semana. "
$ws.Range("G45").Value = "Produz 20 unidades de Ferro toda semana. "
$ws.Range("G46").Value = "Produz 50 unidades de Ferro toda semana. "
$ws.Range("G53").Value = "Produz 20 unidades de Ferro toda semana. "
$ws.Range("G54").Value = "Produz 50 unidades de Ferro toda semana. "
$ws.Range("G61").Value = "Produz 20 unidades de Ferro toda semana. "
$ws.Range("G62").Value = "Produz 50 unidades de Ferro toda semana. "
$ws.Range("G69").Value = "Produz 20 unidades de Ferro toda semana. "
$ws.Range("G70").Value = "Produz 50 unidades de Ferro toda semana. "
$ws.Range("G77").Value = "Produz 20 unidades de Ferro toda semana. "
$ws.Range("G78").Value = "Produz 50 unidades de Ferro toda semana. "
$ws.Range("G86").Value = "Produz 20 unidades de Ferro toda semana. "
$ws.Range("G87").Value = "Produz 50 unidades de Ferro toda semana. "
$ws.Range("G88").Value = "Produz 100 Ferro toda semana."
$ws.Range("G104").Value = "Produz 20 unidades de Ferro toda semana. "
$ws.Range("G105").Value = "Produz 50 unidades de Ferro toda semana. "
$ws.Range("G106").Value = "Produz 100 Ferro toda semana."
$ws.Range("G125").Value = "Produz 20 unidades de Ferro toda semana. "
$ws.Range("G126").Value = "Produz 50 unidades de Ferro toda semana. "
$ws.Range("G127").Value = "Produz 100 Ferro toda semana."
$ws.Range("G142").Value = "Produz 20 unidades de Ferro toda semana. "
$ws.Range("G143").Value = "Produz 50 unidades de Ferro toda semana. "
$ws.Range("G144").Value = "Produz 100 Ferro toda semana."
$ws.Range("G159").Value = "Produz 20 unidades de Ferro toda semana. "
$ws.Range("G160").Value = "Produz 50 unidades de Ferro toda semana. "
$ws.Range("G161").Value = "Produz 100 Ferro toda semana."
$ws.Range("G180").Value = "Produz 20 unidades de Ferro toda semana. "
$ws.Range("G181").Value = "Produz 50 unidades de Ferro toda semana. "
$ws.Range("G182").Value = "Produz 100 Ferro toda semana."
$ws.Range("G210").Value = "O Domínio Ifrit é rico em Jade e Ferro. Produz 30 Pedras de Ferro e 3 Jades toda semana."
$ws.Range("G211").Value = "O Domínio Ifrit é rico em Jade e Ferro. Produz 50 Pedras de Ferro e 5 Jades toda semana."
$ws.Range("G218").Value = "Produz 20 unidades de Ferro toda semana. "
$ws.Range("G219").Value = "Produz 50 unidades de Ferro toda semana. "
$ws.Range("G220").Value = "Produz 100 Ferro toda semana."

# --- Update the sheet view: scroll back to the top (G2) and select G2 ---
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 7
$ws.Range("G2").Select()
